# Tables.xlsx "commit before sector scaffold"
# - Add design_status / cell / sector_id to the shared-string pool (done
#   implicitly by writing them into cells below).
# - Populate new columns E (Cell) and G (Antenna) with field "type" info,
#   shuffle column D (Sector) field names, add two new Antenna rows
#   (design_status / sector_id).
# - Make the "cell" sheet the active sheet/tab, with G12 selected.

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("cell")

# --- Column D ("Sector" field names) updates -----------------------------
# A new "cell" entry is inserted at D3, pushing the old D3/D4 down one row;
# "blocking" and "site_id" move from D5/(new) down to D9/D10.
$ws2.Range("D3").Value = "cell"
$ws2.Range("D4").Value = "serving_area"
$ws2.Range("D5").Value = "morphology"

# D9/D10 are new cells; D10 needs to pick up the same style as D6/D7/D8
# (cellXfs index 12), so copy formats from D6 before setting the value.
$ws2.Range("D9").Value = "blocking"
$ws2.Range("D6").Copy()
$ws2.Range("D10").PasteSpecial(-4122)
$ws2.Range("D10").Value = "site_id"

# --- Column E ("Cell" field types, new) -----------------------------------
$ws2.Range("E2").Value = "string"
$ws2.Range("E3").Value = "integer"
$ws2.Range("E4").Value = "string"
$ws2.Range("E5").Value = "string"
$ws2.Range("E6").Value = "string"
$ws2.Range("E7").Value = "string"
$ws2.Range("E8").Value = "decimal"
$ws2.Range("E9").Value = "string"
$ws2.Range("E10").Value = "integer"

# --- Column F ("Antenna" field names): two new rows -----------------------
$ws2.Range("F11").Value = "design_status"
$ws2.Range("F12").Value = "sector_id"

# --- Column G ("Antenna" field types, new) --------------------------------
$ws2.Range("G2").Value = "string"
$ws2.Range("G3").Value = "string"
$ws2.Range("G4").Value = "string"
$ws2.Range("G5").Value = "decimal"
$ws2.Range("G6").Value = "decimal"
$ws2.Range("G7").Value = "decimal"
$ws2.Range("G8").Value = "decimal"
$ws2.Range("G9").Value = "decimal"
$ws2.Range("G10").Value = "decimal"
$ws2.Range("G11").Value = "string"
$ws2.Range("G12").Value = "integer"

# --- Sheet views / active tab ---------------------------------------------
# "cell" (sheet2) becomes the active/tabSelected sheet with G12 selected;
# "site" (sheet1) loses tabSelected (its own selection, G2, is unchanged).
$ws2.Activate()
$ws2.Range("G12").Select()
